$wb = $excel.ActiveWorkbook

# --- Rename existing sheet, add the new "Login" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Registration"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Login"

# --- Populate the new "Login" sheet ---
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"

$ws2.Range("A2").Value = "tester1@gmail.com"
$ws2.Range("B2").Value = "test123"

$ws2.Range("A3").Value = "tester@gmail.com"
$ws2.Range("B3").Value = "tester123"

$ws2.Range("A4").Value = "tester@gmail.com"
$ws2.Range("B4").Value = "test123"

$ws2.Range("A6").Value = "tester1@gmail.com"
$ws2.Range("B6").Value = "tester123"

$ws2.Range("A7").Value = "demo4@example.com"
$ws2.Range("B7").Value = "test1234"

$ws2.Range("A8").Value = "demo4@example.com"
$ws2.Range("B8").Value = "test1234"

$ws2.Range("A9").Value = "demo4@example.com"
$ws2.Range("B9").Value = "test1234"

# column A width on the Login sheet
$ws2.Columns.Item(1).ColumnWidth = 16.453125

# --- Selections / active sheet ---
$ws1.Range("B8").Select() | Out-Null
$ws2.Range("B5").Select() | Out-Null
